$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.379.76'
$ws.Range("E2").Value = '  +4.80%  '
$ws.Range("D3").Value = '1.812.61'
$ws.Range("E3").Value = '  +5.54%  '
$ws.Range("D4").Value = "'1.004"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").Value = "'343.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.10%  '
$ws.Range("D6").Value = "'1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.19%  '
$ws.Range("D7").Value = "'0.3811"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +3.17%  '
$ws.Range("D8").Value = "'0.3490"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.40%  '
$ws.Range("D9").Value = "'48.91"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.47%  '
$ws.Range("D10").Value = "'1.229"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.39%  '
$ws.Range("D11").Value = "'0.07723"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.41%  '
$ws.Range("D12").Value = "'1.003"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +0.23%  '
$ws.Range("D13").Value = "'22.04"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +9.58%  '
$ws.Range("D14").Value = "'6.593"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +4.33%  '
$ws.Range("D15").Value = '1.812.31'
$ws.Range("E15").Value = '  +5.34%  '
$ws.Range("D16").Value = "'7.201"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.70%  '
$ws.Range("D17").Value = "'0.00001116"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.46%  '
$ws.Range("D18").Value = "'0.06719"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.11%  '
$ws.Range("D19").Value = "'86.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.90%  '
$ws.Range("D20").Value = "'1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("D21").Value = "'17.56"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +6.82%  '
$ws.Range("D22").Value = "'6.550"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +7.47%  '
$ws.Range("D23").Value = "'13.20"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.33%  '
$ws.Range("D24").Value = '27.387.77'
$ws.Range("E24").Value = '  +5.01%  '
$ws.Range("D25").Value = "'2.470"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.18%  '
$ws.Range("D26").Value = "'2.653"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +7.65%  '
$ws.Range("D27").Value = "'21.99"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +14.16%  '
$ws.Range("D28").Value = "'1.461"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +6.66%  '
$ws.Range("D29").Value = "'153.87"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.85%  '
$ws.Range("D30").Value = '2.019.07'
$ws.Range("E30").Value = '  +5.55%  '
$ws.Range("D31").Value = "'135.64"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +4.68%  '
$ws.Range("D32").Value = "'6.295"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.40%  '
$ws.Range("D33").Value = "'4.034"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -2.15%  '
$ws.Range("D34").Value = "'13.85"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +7.09%  '
$ws.Range("D35").Value = "'0.08774"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +2.50%  '
$ws.Range("D36").Value = "'1.688"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.73%  '
$ws.Range("D37").Value = "'5.599"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +4.42%  '
$ws.Range("D38").Value = "'0.6939"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +12.34%  '
$ws.Range("D39").Value = "'0.2267"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +5.79%  '
$ws.Range("D40").Value = "'0.02394"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +3.46%  '
$ws.Range("D41").Value = "'0.06458"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.38%  '
$ws.Range("D42").Value = "'8.895"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +3.31%  '
$ws.Range("D43").Value = "'1.298"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +5.27%  '
$ws.Range("D44").Value = "'14.72"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +1.32%  '
$ws.Range("D45").Value = "'0.6493"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +9.90%  '
$ws.Range("D46").Value = "'1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.13%  '
$ws.Range("D47").Value = "'4.017"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +4.57%  '
$ws.Range("D48").Value = "'2.167"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +7.07%  '
$ws.Range("D49").Value = "'132.55"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.10%  '
$ws.Range("D50").Value = "'0.07321"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.77%  '
$ws.Range("D51").Value = "'80.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +4.09%  '
